# Update Sheets via scheduled runner: refresh currentAveragePrice / Leve profit
# columns (H-N) for a set of leve rows across multiple job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 502002.5
$ws.Range("I113").Value = 502002.5
$ws.Range("K113").Value = 502002.5
$ws.Range("M113").Value = -498748.5

$ws.Range("H134").Value = 41666.668
$ws.Range("J134").Value = 41666.668
$ws.Range("L134").Value = 41666.668
$ws.Range("N134").Value = -51806.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 5333.3335
$ws.Range("J12").Value = 5333.3335
$ws.Range("L12").Value = 5333.3335
$ws.Range("N12").Value = -5679.3335

$ws.Range("H102").Value = 85237.5
$ws.Range("I102").Value = 201618
$ws.Range("J102").Value = 2108.5715
$ws.Range("K102").Value = 201618
$ws.Range("L102").Value = 2108.5715
$ws.Range("M102").Value = -199996
$ws.Range("N102").Value = -5352.5715

$ws.Range("H122").Value = 1947.2
$ws.Range("I122").Value = 1792.8
$ws.Range("K122").Value = 5378.4
$ws.Range("M122").Value = -2928.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H94").Value = 905.44446
$ws.Range("I94").Value = 928.4286
$ws.Range("J94").Value = 825
$ws.Range("K94").Value = 928.4286
$ws.Range("L94").Value = 825
$ws.Range("M94").Value = -477.4286
$ws.Range("N94").Value = -1727

$ws.Range("H105").Value = 78483.38
$ws.Range("I105").Value = 64127.312
$ws.Range("J105").Value = 101453.1
$ws.Range("K105").Value = 64127.312
$ws.Range("L105").Value = 101453.1
$ws.Range("M105").Value = -62380.312
$ws.Range("N105").Value = -104947.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2833.6365
$ws.Range("I58").Value = 2392.75
$ws.Range("J58").Value = 4009.3333
$ws.Range("K58").Value = 2392.75
$ws.Range("L58").Value = 4009.3333
$ws.Range("M58").Value = -2189.75
$ws.Range("N58").Value = -4415.3333

$ws.Range("H136").Value = 2833.6365
$ws.Range("I136").Value = 2392.75
$ws.Range("J136").Value = 4009.3333
$ws.Range("K136").Value = 7178.25
$ws.Range("L136").Value = 12027.9999
$ws.Range("M136").Value = -4628.25
$ws.Range("N136").Value = -17127.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 117878.72
$ws.Range("I70").Value = 131600.5
$ws.Range("J70").Value = 8104.5
$ws.Range("K70").Value = 131600.5
$ws.Range("L70").Value = 8104.5
$ws.Range("M70").Value = -131330.5
$ws.Range("N70").Value = -8644.5

$ws.Range("H73").Value = 117878.72
$ws.Range("I73").Value = 131600.5
$ws.Range("J73").Value = 8104.5
$ws.Range("K73").Value = 131600.5
$ws.Range("L73").Value = 8104.5
$ws.Range("M73").Value = -130664.5
$ws.Range("N73").Value = -9976.5

$ws.Range("H80").Value = 166841500
$ws.Range("J80").Value = 10000
$ws.Range("L80").Value = 10000
$ws.Range("N80").Value = -11996

$ws.Range("H83").Value = 166841500
$ws.Range("J83").Value = 10000
$ws.Range("L83").Value = 50000
$ws.Range("N83").Value = -59984

$ws.Range("H122").Value = 2351.625
$ws.Range("I122").Value = 1839
$ws.Range("K122").Value = 5517
$ws.Range("M122").Value = -3067

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6178.933
$ws.Range("I7").Value = 7025.5
$ws.Range("J7").Value = 5211.4287
$ws.Range("K7").Value = 7025.5
$ws.Range("L7").Value = 5211.4287
$ws.Range("M7").Value = -6913.5
$ws.Range("N7").Value = -5435.4287

$ws.Range("H40").Value = 48350.547
$ws.Range("I40").Value = 86759.336
$ws.Range("K40").Value = 86759.336
$ws.Range("M40").Value = -86623.336

$ws.Range("H68").Value = 3306.4
$ws.Range("I68").Value = 1555.2727
$ws.Range("J68").Value = 5446.6665
$ws.Range("K68").Value = 1555.2727
$ws.Range("L68").Value = 5446.6665
$ws.Range("M68").Value = -806.2727
$ws.Range("N68").Value = -6944.6665

$ws.Range("H71").Value = 3306.4
$ws.Range("I71").Value = 1555.2727
$ws.Range("J71").Value = 5446.6665
$ws.Range("K71").Value = 7776.363499999999
$ws.Range("L71").Value = 27233.3325
$ws.Range("M71").Value = -4032.363499999999
$ws.Range("N71").Value = -34721.3325

$ws.Range("H93").Value = 1576.8235
$ws.Range("I93").Value = 1496.9565
$ws.Range("J93").Value = 1743.8182
$ws.Range("K93").Value = 1496.9565
$ws.Range("L93").Value = 1743.8182
$ws.Range("M93").Value = -248.9565
$ws.Range("N93").Value = -4239.8182

$ws.Range("H100").Value = 1943.4546
$ws.Range("I100").Value = 1630
$ws.Range("J100").Value = 2204.6667
$ws.Range("K100").Value = 1630
$ws.Range("L100").Value = 2204.6667
$ws.Range("M100").Value = -1089
$ws.Range("N100").Value = -3286.6667

$ws.Range("H126").Value = 6178.933
$ws.Range("I126").Value = 7025.5
$ws.Range("J126").Value = 5211.4287
$ws.Range("K126").Value = 21076.5
$ws.Range("L126").Value = 15634.2861
$ws.Range("M126").Value = -18606.5
$ws.Range("N126").Value = -20574.2861

$ws.Range("H132").Value = 4290.5415
$ws.Range("I132").Value = 4169.0625
$ws.Range("J132").Value = 4533.5
$ws.Range("K132").Value = 12507.1875
$ws.Range("L132").Value = 13600.5
$ws.Range("M132").Value = -9977.1875
$ws.Range("N132").Value = -18660.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2326.158
$ws.Range("I122").Value = 1292
$ws.Range("J122").Value = 2929.4167
$ws.Range("K122").Value = 3876
$ws.Range("L122").Value = 8788.250100000001
$ws.Range("M122").Value = -1426
$ws.Range("N122").Value = -13688.2501

$ws.Range("H126").Value = 1926.1428
$ws.Range("I126").Value = 1821
$ws.Range("J126").Value = 2066.3333
$ws.Range("K126").Value = 5463
$ws.Range("L126").Value = 6198.999899999999
$ws.Range("M126").Value = -2993
$ws.Range("N126").Value = -11138.9999

$ws.Range("H132").Value = 2434.0222
$ws.Range("I132").Value = 2640.3225
$ws.Range("J132").Value = 1977.2142
$ws.Range("K132").Value = 7920.967500000001
$ws.Range("L132").Value = 5931.642599999999
$ws.Range("M132").Value = -5390.967500000001
$ws.Range("N132").Value = -10991.6426

$ws.Range("H136").Value = 1595.3396
$ws.Range("I136").Value = 626.7917
$ws.Range("J136").Value = 2396.8965
$ws.Range("K136").Value = 1880.3751
$ws.Range("L136").Value = 7190.689499999999
$ws.Range("M136").Value = 669.6249
$ws.Range("N136").Value = -12290.6895

